# Applies the documented changes to "Using JTAG And Hardware Manager.docx":
#  1. "Yes" -> "OK" (bold confirmation answer)
#  2. "saturnprom" -> "saturnfallback" (file name update)
#  3. Insert a new list item "Click OK to begin" before the
#     "The config prom is programmed ..." list item.

$d = $word.ActiveDocument

# --- 1. Replace the bold "Yes" answer with "OK" -----------------------
$rng1 = $d.Content.Duplicate
$rng1.Find.ClearFormatting()
$found1 = $rng1.Find.Execute("Yes", $true, $false, $false, $false, $false, $true, 1, $false, "OK", 2)
Write-Host "Step 1 (Yes -> OK): $found1"

# --- 2. Update the file name referenced in the instructions ------------
$rng2 = $d.Content.Duplicate
$rng2.Find.ClearFormatting()
$found2 = $rng2.Find.Execute("saturnprom", $true, $false, $false, $false, $false, $true, 1, $false, "saturnfallback", 2)
Write-Host "Step 2 (saturnprom -> saturnfallback): $found2"

# --- 3. Insert a new bullet before "The config prom is programmed..." --
$rng3 = $d.Content.Duplicate
$found3 = $rng3.Find.Execute("The config prom is programmed", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Step 3 locate target paragraph: $found3"

if ($found3) {
    $targetPara = $rng3.Paragraphs(1)
    $paraStart = $targetPara.Range.Start
    $insertPoint = $d.Range($paraStart, $paraStart)
    $insertPoint.InsertBefore("Click OK to begin`r")
    Write-Host "Step 3: inserted new list paragraph"
}

Write-Host "Done"
